$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number/date need to be forced back to
# Text format around the assignment, otherwise Excel auto-converts the inline
# string into a numeric cell. Style is reset to "Normal" afterwards so no residual
# number-format attribute is left on the cell (matches original formatting).

$ws.Range('D2').Value = '67.832.25'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '3.266.07'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.05%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.601'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('E9').Value = '  -2.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('E11').Value = '  -3.61%  '
$ws.Range('D12').Value = '3.831.91'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.35'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.14%  '
$ws.Range('D15').Value = '67.811.49'
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('E16').Value = '  -1.90%  '
$ws.Range('D17').Value = '3.266.92'
$ws.Range('E17').Value = '  +2.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.69'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.39'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '402.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.54'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.86%  '
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.508'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.45%  '
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.188'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.46'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  -1.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.65'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.46'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.03%  '
$ws.Range('B32').Value = 'USDe'
$ws.Range('C32').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('B33').Value = 'Aptos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.88'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.58%  '
$ws.Range('E34').Value = '  -2.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '164.47'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E36').Value = '  -3.71%  '
$ws.Range('E37').Value = '  -0.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.11'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.799'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.57%  '
$ws.Range('E40').Value = '  -2.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.33'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Value = '2.679.87'
$ws.Range('E42').Value = '  +2.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.82'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0677'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.43'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '334.61'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.29%  '
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('E48').Value = '  -2.37%  '
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('E51').Value = '  -1.40%  '
